$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the existing row 169 (shifts old 169-173 down to 172-176)
$ws.Rows("169:171").Insert()

# Row 169: new weekly entry - Extra
$ws.Cells.Item(169, 1).Value = 7
$ws.Cells.Item(169, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(169, 3).Value = "Ñuble"
$ws.Cells.Item(169, 4).Value = 44610
$ws.Cells.Item(169, 5).Value = 16
$ws.Cells.Item(169, 6).Value = 100112028
$ws.Cells.Item(169, 7).Value = "Sandia"
$ws.Cells.Item(169, 8).Value = "Sin especificar"
$ws.Cells.Item(169, 9).Value = "Extra"
$ws.Cells.Item(169, 10).Value = 500
$ws.Cells.Item(169, 11).Value = 2500
$ws.Cells.Item(169, 12).Value = 2500
$ws.Cells.Item(169, 13).Value = 2500
$ws.Cells.Item(169, 14).Value = "$/unidad"
$ws.Cells.Item(169, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(169, 16).Value = 2500
$ws.Cells.Item(169, 17).Value = 1
$ws.Cells.Item(169, 18).Value = "Hortaliza"

# Row 170: new weekly entry - Primera
$ws.Cells.Item(170, 1).Value = 7
$ws.Cells.Item(170, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(170, 3).Value = "Ñuble"
$ws.Cells.Item(170, 4).Value = 44610
$ws.Cells.Item(170, 5).Value = 16
$ws.Cells.Item(170, 6).Value = 100112028
$ws.Cells.Item(170, 7).Value = "Sandia"
$ws.Cells.Item(170, 8).Value = "Sin especificar"
$ws.Cells.Item(170, 9).Value = "Primera"
$ws.Cells.Item(170, 10).Value = 600
$ws.Cells.Item(170, 11).Value = 2000
$ws.Cells.Item(170, 12).Value = 2200
$ws.Cells.Item(170, 13).Value = 2100
$ws.Cells.Item(170, 14).Value = "$/unidad"
$ws.Cells.Item(170, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(170, 16).Value = 2100
$ws.Cells.Item(170, 17).Value = 1
$ws.Cells.Item(170, 18).Value = "Hortaliza"

# Row 171: new weekly entry - Segunda
$ws.Cells.Item(171, 1).Value = 7
$ws.Cells.Item(171, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(171, 3).Value = "Ñuble"
$ws.Cells.Item(171, 4).Value = 44610
$ws.Cells.Item(171, 5).Value = 16
$ws.Cells.Item(171, 6).Value = 100112028
$ws.Cells.Item(171, 7).Value = "Sandia"
$ws.Cells.Item(171, 8).Value = "Sin especificar"
$ws.Cells.Item(171, 9).Value = "Segunda"
$ws.Cells.Item(171, 10).Value = 600
$ws.Cells.Item(171, 11).Value = 1500
$ws.Cells.Item(171, 12).Value = 1700
$ws.Cells.Item(171, 13).Value = 1600
$ws.Cells.Item(171, 14).Value = "$/unidad"
$ws.Cells.Item(171, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(171, 16).Value = 1600
$ws.Cells.Item(171, 17).Value = 1
$ws.Cells.Item(171, 18).Value = "Hortaliza"
